$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.839.54"
$ws.Range("E2").Value = "  +6.84%  "
$ws.Range("D3").Value = "3.021.81"
$ws.Range("E3").Value = "  +4.36%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.77"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.60"
$ws.Range("E6").Value = "  +8.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.017.37"
$ws.Range("E8").Value = "  +4.38%  "
$ws.Range("E9").Value = "  +2.82%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("E12").Value = "  +5.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +6.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.56"
$ws.Range("E14").Value = "  +8.40%  "
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("D16").Value = "65.945.39"
$ws.Range("E16").Value = "  +7.01%  "
$ws.Range("D17").Value = "3.521.93"
$ws.Range("E17").Value = "  +4.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.00"
$ws.Range("E18").Value = "  +6.88%  "
$ws.Range("D19").Value = "3.016.88"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.34"
$ws.Range("E20").Value = "  +6.83%  "
$ws.Range("E21").Value = "  +5.80%  "
$ws.Range("E22").Value = "  +4.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.40"
$ws.Range("E23").Value = "  +8.47%  "
$ws.Range("E24").Value = "  +3.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.61"
$ws.Range("E25").Value = "  +5.22%  "
$ws.Range("E26").Value = "  +12.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  +9.08%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.99"
$ws.Range("E29").Value = "  +14.31%  "
$ws.Range("E30").Value = "  +18.86%  "
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  +4.45%  "
$ws.Range("E33").Value = "  +5.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.08"
$ws.Range("E34").Value = "  +5.93%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +4.08%  "
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.81"
$ws.Range("E37").Value = "  +7.85%  "
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  +13.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.06"
$ws.Range("E39").Value = "  +8.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.41"
$ws.Range("E40").Value = "  +1.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "45.30"
$ws.Range("E41").Value = "  +15.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.123"
$ws.Range("E42").Value = "  +7.63%  "
$ws.Range("E43").Value = "  +13.60%  "
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "393.92"
$ws.Range("E45").Value = "  +13.55%  "
$ws.Range("D46").Value = "2.804.64"
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0354"
$ws.Range("E47").Value = "  +5.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.96"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.87"
$ws.Range("E50").Value = "  +10.62%  "
$ws.Range("E51").Value = "  +4.19%  "

# Restore default General format/style on cells that were forced to Text
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
